$d = $word.ActiveDocument

# Locate the paragraph that starts the footer block ("Ver no Jupiter ...")
# and the paragraph that ends it ("... Creative Commons Attribution").
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($startPara -eq $null) -and ($t -like "*Ver no Jupiter*")) {
        $startPara = $i
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $i
    }
}

if (($startPara -ne $null) -and ($endPara -ne $null)) {
    # Also swallow the now-redundant blank paragraph that trails the footer
    # block (mirrors the blank paragraph that still precedes it), so the
    # document collapses back to a single blank line before the page break.
    $removeEnd = $endPara
    $next = $d.Paragraphs.Item($endPara + 1)
    if ($next.Range.Text.Trim().Length -eq 0) {
        $removeEnd = $endPara + 1
    }

    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($removeEnd).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
